# Replace the 15 lattice-multiplication practice problems in the single
# 5x3 table with a new set of problems/partial-product digits, while
# keeping the table shape, run formatting (sz 32) and the line-break
# layout ("A x B" / "  d    d" / "  ----" / "d|    |" / "d|    |")
# identical.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11   # manual line break -> <w:br/>, same as the existing runs

$t.Cell(1, 1).Range.Text = "35 x 52" + $nl + "  5    2" + $nl + "  ----" + $nl + "3|    |" + $nl + "5|    |"
$t.Cell(1, 2).Range.Text = "58 x 34" + $nl + "  3    4" + $nl + "  ----" + $nl + "5|    |" + $nl + "8|    |"
$t.Cell(1, 3).Range.Text = "51 x 91" + $nl + "  9    1" + $nl + "  ----" + $nl + "5|    |" + $nl + "1|    |"

$t.Cell(2, 1).Range.Text = "93 x 25" + $nl + "  2    5" + $nl + "  ----" + $nl + "9|    |" + $nl + "3|    |"
$t.Cell(2, 2).Range.Text = "48 x 91" + $nl + "  9    1" + $nl + "  ----" + $nl + "4|    |" + $nl + "8|    |"
$t.Cell(2, 3).Range.Text = "91 x 90" + $nl + "  9    0" + $nl + "  ----" + $nl + "9|    |" + $nl + "1|    |"

$t.Cell(3, 1).Range.Text = "16 x 16" + $nl + "  1    6" + $nl + "  ----" + $nl + "1|    |" + $nl + "6|    |"
$t.Cell(3, 2).Range.Text = "28 x 22" + $nl + "  2    2" + $nl + "  ----" + $nl + "2|    |" + $nl + "8|    |"
$t.Cell(3, 3).Range.Text = "43 x 24" + $nl + "  2    4" + $nl + "  ----" + $nl + "4|    |" + $nl + "3|    |"

$t.Cell(4, 1).Range.Text = "99 x 86" + $nl + "  8    6" + $nl + "  ----" + $nl + "9|    |" + $nl + "9|    |"
$t.Cell(4, 2).Range.Text = "35 x 43" + $nl + "  4    3" + $nl + "  ----" + $nl + "3|    |" + $nl + "5|    |"
$t.Cell(4, 3).Range.Text = "86 x 89" + $nl + "  8    9" + $nl + "  ----" + $nl + "8|    |" + $nl + "6|    |"

$t.Cell(5, 1).Range.Text = "80 x 20" + $nl + "  2    0" + $nl + "  ----" + $nl + "8|    |" + $nl + "0|    |"
$t.Cell(5, 2).Range.Text = "92 x 97" + $nl + "  9    7" + $nl + "  ----" + $nl + "9|    |" + $nl + "2|    |"
$t.Cell(5, 3).Range.Text = "93 x 82" + $nl + "  8    2" + $nl + "  ----" + $nl + "9|    |" + $nl + "3|    |"
